$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column CP (94th column), shifting CP:OZ -> CT:PD
$ws.Range("CP1:CS1").EntireColumn.Insert()

# Force explicit (blank) cell records for the inserted columns in the data
# rows, matching how the rest of the sheet pads every row out to the full
# column range. Touching the (already-default) border avoids minting a new
# style index.
$ws.Range("CP2:CS5").Borders.LineStyle = -4142

# Set header labels for the newly inserted columns
$ws.Range("CP1").Value = "MemberLanguageInfoCodeQualifier"
$ws.Range("CQ1").Value = "MemberLanguageInfoCode"
$ws.Range("CR1").Value = "MemberLanguageInfoLanguageDescription"
$ws.Range("CS1").Value = "MemberLanguageInfoLanguageUseIndicator"

# Update the Id value for the data rows
$ws.Range("A2").Value = "685d5721e29351212dc46be9"
$ws.Range("A3").Value = "685d5721e29351212dc46be9"
$ws.Range("A4").Value = "685d5721e29351212dc46bea"
$ws.Range("A5").Value = "685d5721e29351212dc46beb"
